$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.368.69'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = '3.112.53'
$ws.Range("E3").Value = '  -1.85%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '241.48'
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").Value = '615.16'
$ws.Range("E6").Value = '  -1.24%  '

$ws.Range("E7").Value = '  -3.88%  '

$ws.Range("D8").Value = '0.395'
$ws.Range("E8").Value = '  +5.44%  '

$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("D10").Value = '3.108.69'
$ws.Range("E10").Value = '  -1.91%  '

$ws.Range("E11").Value = '  -2.84%  '

$ws.Range("E12").Value = '  -1.19%  '

$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  +0.50%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '34.46'
$ws.Range("E14").Value = '  -3.11%  '

$ws.Range("B15").Value = 'Toncoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D15").Value = '5.53'
$ws.Range("E15").Value = '  +0.21%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '92.027.92'
$ws.Range("E16").Value = '  +0.58%  '

$ws.Range("D17").Value = '3.685.04'
$ws.Range("E17").Value = '  -1.32%  '

$ws.Range("D18").Value = '3.094.71'
$ws.Range("E18").Value = '  -1.51%  '

$ws.Range("E19").Value = '  -2.14%  '

$ws.Range("D20").Value = '14.75'
$ws.Range("E20").Value = '  -7.07%  '

$ws.Range("D21").Value = '5.83'
$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").Value = '9.46'
$ws.Range("E22").Value = '  +1.14%  '

$ws.Range("D23").Value = '448.28'
$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("E24").Value = '  -5.10%  '

$ws.Range("D25").Value = '5.62'
$ws.Range("E25").Value = '  -7.41%  '

$ws.Range("D26").Value = '87.03'
$ws.Range("E26").Value = '  -3.54%  '

$ws.Range("D27").Value = '11.75'
$ws.Range("E27").Value = '  -2.55%  '

$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("D30").Value = '0.139'
$ws.Range("E30").Value = '  -5.52%  '

$ws.Range("D31").Value = '0.235'
$ws.Range("E31").Value = '  -1.51%  '

$ws.Range("D32").Value = '0.169'
$ws.Range("E32").Value = '  -3.06%  '

$ws.Range("D33").Value = '9.26'
$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  +23.87%  '

$ws.Range("D35").Value = '8.10'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("D36").Value = '0.165'
$ws.Range("E36").Value = '  -2.35%  '

$ws.Range("D37").Value = '4.24'
$ws.Range("E37").Value = '  -0.49%  '

$ws.Range("D38").Value = '26.24'
$ws.Range("E38").Value = '  -2.04%  '

$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -0.60%  '

$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("D41").Value = '481.77'
$ws.Range("E41").Value = '  -6.44%  '

$ws.Range("D42").Value = '0.437'
$ws.Range("E42").Value = '  -1.20%  '

$ws.Range("D43").Value = '3.46'
$ws.Range("E43").Value = '  -1.42%  '

$ws.Range("D44").Value = '23.05'
$ws.Range("E44").Value = '  +3.91%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = '158.17'
$ws.Range("E46").Value = '  +1.79%  '

$ws.Range("D47").Value = '1.91'
$ws.Range("E47").Value = '  -1.82%  '

$ws.Range("D48").Value = '0.695'
$ws.Range("E48").Value = '  -3.10%  '

$ws.Range("E49").Value = '  -0.71%  '

$ws.Range("D50").Value = '0.0337'
$ws.Range("E50").Value = '  +2.87%  '

$ws.Range("E51").Value = '  -0.75%  '
